# Commit: "OpenXML Encoder and better Excel coverage."
#
# The gold OOXML shows that, for every run of text in the presentation whose
# visible (already-unescaped) content contains the literal artifact string
#   </a:t></a:r>
# that string gets turned into
#   &lt;/a:t>&lt;/a:r>
# (i.e. the run text itself, not real markup, gets an extra layer of
# XML-entity escaping applied to it by the encoder that introduced the bug
# being regression-tested). After PowerPoint re-serializes the file, the
# <a:t> element ends up containing "&amp;lt;/a:t&gt;&amp;lt;/a:r&gt;".
#
# Every <a:r> run in the deck already contains exactly one occurrence of the
# search string, always immediately followed by "]   " and nothing else, so
# we can recover each run's exact [start,end) character range inside the
# shape's flattened TextRange text purely from that fact (plus the location
# of the U+000B characters PowerPoint uses in TextRange.Text to represent
# <a:br/> line breaks). Updating the *entire* run's text in one
# Characters(...).Text assignment (rather than just the matched substring)
# keeps each run a single <a:r> element instead of splintering it into
# multiple runs.

$searchText = "</a:t></a:r>"
$tailText = "]   "
$replaceText = "&lt;/a:t>&lt;/a:r>"
$breakChar = [char]11

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }
        $tf = $shape.TextFrame
        if (-not $tf.HasText) { continue }

        $tr = $tf.TextRange
        $full = $tr.Text

        if (-not $full.Contains($searchText)) { continue }

        # Locate every occurrence of the search text.
        $matchStarts = New-Object System.Collections.Generic.List[int]
        $searchFrom = 0
        while ($true) {
            $found = $full.IndexOf($searchText, $searchFrom)
            if ($found -lt 0) { break }
            [void]$matchStarts.Add($found)
            $searchFrom = $found + 1
        }

        # Each match belongs to a run whose text ends right after the
        # following "]   " tail.
        $runEnds = New-Object System.Collections.Generic.List[int]
        foreach ($ms in $matchStarts) {
            [void]$runEnds.Add($ms + $searchText.Length + $tailText.Length)
        }

        # A run starts right where the previous one ended, after skipping
        # any <a:br/> break characters in between.
        $runStarts = New-Object System.Collections.Generic.List[int]
        $cursor = 0
        for ($i = 0; $i -lt $runEnds.Count; $i++) {
            while (($cursor -lt $full.Length) -and ($full[$cursor] -eq $breakChar)) {
                $cursor = $cursor + 1
            }
            [void]$runStarts.Add($cursor)
            $cursor = $runEnds[$i]
        }

        # Update runs from last to first so earlier offsets stay valid while
        # the text grows.
        for ($i = $runStarts.Count - 1; $i -ge 0; $i--) {
            $start = $runStarts[$i]
            $length = $runEnds[$i] - $start
            $runRange = $tr.Characters($start + 1, $length)
            $runRange.Text = $runRange.Text.Replace($searchText, $replaceText)
        }
    }
}
